$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new content to row 3 (B3, C3)
$ws.Range("B3").Value = "due to a coding mistake of writing phase inside raise.objection()"
$ws.Range("C3").Value = "write this in raise.objection() "

# C3: vertical centered, no wrap (create this style first -> xf index 4)
$ws.Range("C3").VerticalAlignment = -4108
$ws.Range("C3").WrapText = $false

# B3: vertical centered, wrap text (create this style second -> xf index 5)
$ws.Range("B3").VerticalAlignment = -4108
$ws.Range("B3").WrapText = $true

# Let Excel autofit the row height for the wrapped text
$ws.Rows.Item(3).AutoFit() | Out-Null
Write-Host "RowHeight after autofit:" $ws.Rows.Item(3).RowHeight

# Update the active cell selection to B3
$ws.Range("B3").Select()
